$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Delete rows 12 through 33 (shrinks the used range from A1:J33 down to A1:J11)
$ws.Range("A12:J33").EntireRow.Delete()

# Restore the view/selection state seen in the target workbook
$ws.Application.ActiveWindow.ScrollRow = 2
$ws.Range("J16").Select()
